$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 522, shifting rows 522:546 down to 523:547.
$ws.Rows.Item(522).Insert()

# Populate the newly inserted row 522 with the new data record.
$ws.Range("A522").Value = 4
$ws.Range("B522").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C522").Value = "Los Lagos"
$ws.Range("D522").Value = 44931
$ws.Range("E522").Value = 10
$ws.Range("F522").Value = "Fruta"
$ws.Range("G522").Value = 100103
$ws.Range("H522").Value = "Frutos de hueso (carozo)"
$ws.Range("I522").Value = 100103006
$ws.Range("J522").Value = "Nectarín"
$ws.Range("K522").Value = "Super Queen"
$ws.Range("L522").Value = "Especial"
$ws.Range("M522").Value = 400
$ws.Range("N522").Value = 23000
$ws.Range("O522").Value = 23000
$ws.Range("P522").Value = 23000
$ws.Range("Q522").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R522").Value = "Región de O'Higgins"
$ws.Range("S522").Value = 1643
$ws.Range("T522").Value = 14
